$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "63.805.50"
Set-TextCell "E2" "  -4.83%  "

# Row 3
Set-TextCell "D3" "3.001.03"
Set-TextCell "E3" "  -6.74%  "

# Row 4
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  +0.01%  "

# Row 5
Set-TextCell "D5" "551.11"
Set-TextCell "E5" "  -7.69%  "

# Row 6
Set-TextCell "D6" "137.37"
Set-TextCell "E6" "  -8.65%  "

# Row 7
Set-TextCell "D7" "0.995"
Set-TextCell "E7" "  -0.60%  "

# Row 8
Set-TextCell "D8" "2.979.48"
Set-TextCell "E8" "  -7.12%  "

# Row 9
Set-TextCell "D9" "0.474"
Set-TextCell "E9" "  -12.69%  "

# Row 10
Set-TextCell "D10" "0.151"
Set-TextCell "E10" "  -13.35%  "

# Row 11
Set-TextCell "D11" "6.14"
Set-TextCell "E11" "  -6.16%  "

# Row 12
Set-TextCell "D12" "0.444"
Set-TextCell "E12" "  -10.73%  "

# Row 13
Set-TextCell "D13" "33.61"
Set-TextCell "E13" "  -13.54%  "

# Row 14
Set-TextCell "D14" "0.0000211"
Set-TextCell "E14" "  -13.92%  "

# Row 15
Set-TextCell "D15" "3.467.93"
Set-TextCell "E15" "  -7.15%  "

# Row 16
Set-TextCell "D16" "63.682.07"
Set-TextCell "E16" "  -4.98%  "

# Row 17
Set-TextCell "E17" "  -4.58%  "

# Row 18
Set-TextCell "D18" "2.985.51"
Set-TextCell "E18" "  -7.49%  "

# Row 19
Set-TextCell "D19" "479.55"
Set-TextCell "E19" "  -9.58%  "

# Row 20
Set-TextCell "D20" "6.40"
Set-TextCell "E20" "  -10.70%  "

# Row 21
Set-TextCell "D21" "13.13"
Set-TextCell "E21" "  -11.75%  "

# Row 22
Set-TextCell "D22" "0.649"
Set-TextCell "E22" "  -14.50%  "

# Row 23
Set-TextCell "D23" "6.73"
Set-TextCell "E23" "  -14.84%  "

# Row 24
Set-TextCell "D24" "12.27"
Set-TextCell "E24" "  -11.09%  "

# Row 25
Set-TextCell "D25" "76.94"
Set-TextCell "E25" "  -9.97%  "

# Row 26
Set-TextCell "E26" "  +0.04%  "

# Row 27
Set-TextCell "D27" "2.67"
Set-TextCell "E27" "  -16.12%  "

# Row 28
Set-TextCell "D28" "7.43"
Set-TextCell "E28" "  -8.49%  "

# Row 29
Set-TextCell "D29" "1.96"
Set-TextCell "E29" "  -10.53%  "

# Row 30
Set-TextCell "D30" "25.35"
Set-TextCell "E30" "  -13.11%  "

# Row 31
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D31" "1.11"
Set-TextCell "E31" "  -2.10%  "

# Row 32
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D32" "2.46"
Set-TextCell "E32" "  -6.91%  "

# Row 33
Set-TextCell "D33" "1.00"
Set-TextCell "E33" "  -0.02%  "

# Row 34
Set-TextCell "D34" "498.32"
Set-TextCell "E34" "  -8.20%  "

# Row 35
Set-TextCell "D35" "51.55"
Set-TextCell "E35" "  -3.22%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D36" "5.04"
Set-TextCell "E36" "  -11.46%  "

# Row 37
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D37" "5.68"
Set-TextCell "E37" "  -12.61%  "

# Row 38
Set-TextCell "D38" "0.0398"
Set-TextCell "E38" "  -6.45%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D39" "0.117"
Set-TextCell "E39" "  -7.05%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D40" "0.0766"
Set-TextCell "E40" "  -11.24%  "

# Row 41
Set-TextCell "D41" "8.03"
Set-TextCell "E41" "  -14.03%  "

# Row 42
Set-TextCell "D42" "2.767.27"
Set-TextCell "E42" "  -5.14%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D43" "2.42"
Set-TextCell "E43" "  -7.50%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D44" "0.998"
Set-TextCell "E44" "  -0.19%  "

# Row 45
Set-TextCell "D45" "0.231"
Set-TextCell "E45" "  -12.17%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D46" "1.95"
Set-TextCell "E46" "  -7.55%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D47" "114.86"
Set-TextCell "E47" "  -6.65%  "

# Row 48
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D48" "0.0₃0501"
Set-TextCell "E48" "  -14.10%  "

# Row 49
Set-TextCell "D49" "0.103"
Set-TextCell "E49" "  -9.48%  "

# Row 50
Set-TextCell "D50" "22.58"
Set-TextCell "E50" "  -15.04%  "

# Row 51
Set-TextCell "D51" "1.97"
Set-TextCell "E51" "  -18.29%  "
